# This script updates the "Name of Algo" result data (KNN imputation output)
# by writing the corrected numeric values into columns C and D for the
# affected rows, matching the canonical OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -11.973
$ws.Range("C6").Value = -11.864
$ws.Range("C7").Value = -12.793
$ws.Range("D7").Value = -7.673999999999999
$ws.Range("C8").Value = -12.606
$ws.Range("D11").Value = -7.133000000000001
$ws.Range("D12").Value = -7.514999999999999
$ws.Range("D15").Value = -8.392999999999999
$ws.Range("C16").Value = -13.312
$ws.Range("C20").Value = -12.241
$ws.Range("D20").Value = -7.994
$ws.Range("C21").Value = -12.513
$ws.Range("D21").Value = -8.051999999999998
$ws.Range("D22").Value = -7.536000000000001
$ws.Range("D23").Value = -7.997
$ws.Range("C28").Value = -12.848
$ws.Range("C29").Value = -12.035
$ws.Range("D29").Value = -7.441000000000001
$ws.Range("C30").Value = -12.559
$ws.Range("C32").Value = -12.95
$ws.Range("D34").Value = -7.936999999999999
$ws.Range("C40").Value = -12.445
$ws.Range("D42").Value = -8.107000000000001
$ws.Range("D43").Value = -7.903999999999999
$ws.Range("D44").Value = -7.581
$ws.Range("D45").Value = -7.517
$ws.Range("C46").Value = -13.549
$ws.Range("D46").Value = -8.517000000000001
$ws.Range("D50").Value = -8.248999999999999
$ws.Range("C51").Value = -12.199
$ws.Range("D51").Value = -7.598999999999999
$ws.Range("C52").Value = -11.614
$ws.Range("C57").Value = -13.782
$ws.Range("D57").Value = -8.17
$ws.Range("C59").Value = -12.727
$ws.Range("C62").Value = -13.737
$ws.Range("D65").Value = -7.556999999999999
$ws.Range("C66").Value = -11.418
$ws.Range("D66").Value = -7.683000000000002
$ws.Range("D67").Value = -6.879
$ws.Range("C73").Value = -12.908
$ws.Range("C74").Value = -12.144
$ws.Range("C77").Value = -12.566
$ws.Range("D79").Value = -7.465000000000001
$ws.Range("D84").Value = -8.453000000000001
$ws.Range("D87").Value = -8.068000000000001
$ws.Range("C92").Value = -11.037
$ws.Range("D92").Value = -6.772
$ws.Range("D97").Value = -8.500000000000002
$ws.Range("C100").Value = -12.706
